$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "42.609.61"
Set-TextCell $ws.Range("E2") "  +1.33%  "

Set-TextCell $ws.Range("D3") "2.293.78"
Set-TextCell $ws.Range("E3") "  -0.39%  "

Set-TextCell $ws.Range("E4") "  +0.00%  "

Set-TextCell $ws.Range("D5") "322.71"
Set-TextCell $ws.Range("E5") "  +1.29%  "

Set-TextCell $ws.Range("D6") "103.74"
Set-TextCell $ws.Range("E6") "  -0.87%  "

Set-TextCell $ws.Range("D7") "0.630"
Set-TextCell $ws.Range("E7") "  -0.08%  "

Set-TextCell $ws.Range("E8") "  +0.12%  "

Set-TextCell $ws.Range("D9") "0.610"
Set-TextCell $ws.Range("E9") "  +0.00%  "

Set-TextCell $ws.Range("D10") "40.02"
Set-TextCell $ws.Range("E10") "  +0.84%  "

Set-TextCell $ws.Range("D11") "0.0909"
Set-TextCell $ws.Range("E11") "  -0.51%  "

Set-TextCell $ws.Range("D12") "8.41"
Set-TextCell $ws.Range("E12") "  -0.28%  "

Set-TextCell $ws.Range("D13") "0.107"
Set-TextCell $ws.Range("E13") "  -0.23%  "

Set-TextCell $ws.Range("D14") "0.973"
Set-TextCell $ws.Range("E14") "  -0.66%  "

Set-TextCell $ws.Range("D15") "15.23"
Set-TextCell $ws.Range("E15") "  -1.41%  "

Set-TextCell $ws.Range("D16") "2.638.09"
Set-TextCell $ws.Range("E16") "  -0.63%  "

Set-TextCell $ws.Range("D17") "2.290.84"
Set-TextCell $ws.Range("E17") "  -0.77%  "

Set-TextCell $ws.Range("D18") "42.542.26"
Set-TextCell $ws.Range("E18") "  +0.88%  "

Set-TextCell $ws.Range("D19") "7.42"
Set-TextCell $ws.Range("E19") "  -3.89%  "

Set-TextCell $ws.Range("D20") "0.0000106"

Set-TextCell $ws.Range("D21") "13.66"
Set-TextCell $ws.Range("E21") "  +35.67%  "

Set-TextCell $ws.Range("B22") "Litecoin"
Set-TextCell $ws.Range("C22") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D22") "73.33"
Set-TextCell $ws.Range("E22") "  -0.69%  "

Set-TextCell $ws.Range("B23") "PancakeSwap"
Set-TextCell $ws.Range("C23") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D23") "3.59"
Set-TextCell $ws.Range("E23") "  -1.84%  "

Set-TextCell $ws.Range("D24") "269.57"
Set-TextCell $ws.Range("E24") "  -7.08%  "

Set-TextCell $ws.Range("E25") "  -2.55%  "

Set-TextCell $ws.Range("E26") "  -0.41%  "

Set-TextCell $ws.Range("D27") "10.92"
Set-TextCell $ws.Range("E27") "  -0.47%  "

Set-TextCell $ws.Range("E28") "  +2.79%  "

Set-TextCell $ws.Range("D29") "22.54"
Set-TextCell $ws.Range("E29") "  -3.51%  "

Set-TextCell $ws.Range("D30") "38.21"
Set-TextCell $ws.Range("E30") "  +7.80%  "

Set-TextCell $ws.Range("D31") "165.11"
Set-TextCell $ws.Range("E31") "  +0.51%  "

Set-TextCell $ws.Range("D32") "6.17"
Set-TextCell $ws.Range("E32") "  +4.38%  "

Set-TextCell $ws.Range("D33") "0.0882"
Set-TextCell $ws.Range("E33") "  -0.65%  "

Set-TextCell $ws.Range("E34") "  +0.42%  "

Set-TextCell $ws.Range("E35") "  -1.72%  "

Set-TextCell $ws.Range("E36") "  -14.15%  "

Set-TextCell $ws.Range("D37") "4.63"
Set-TextCell $ws.Range("E37") "  -0.52%  "

Set-TextCell $ws.Range("D38") "0.0356"
Set-TextCell $ws.Range("E38") "  +0.72%  "

Set-TextCell $ws.Range("D39") "3.71"
Set-TextCell $ws.Range("E39") "  +1.99%  "

Set-TextCell $ws.Range("E40") "  -5.11%  "

Set-TextCell $ws.Range("E41") "  +3.37%  "

Set-TextCell $ws.Range("D42") "69.58"
Set-TextCell $ws.Range("E42") "  -1.96%  "

Set-TextCell $ws.Range("E43") "  -0.11%  "

Set-TextCell $ws.Range("D44") "0.226"
Set-TextCell $ws.Range("E44") "  -1.01%  "

Set-TextCell $ws.Range("D45") "92.91"
Set-TextCell $ws.Range("E45") "  -10.16%  "

Set-TextCell $ws.Range("D46") "12.35"
Set-TextCell $ws.Range("E46") "  +1.64%  "

Set-TextCell $ws.Range("D47") "81.79"
Set-TextCell $ws.Range("E47") "  +4.99%  "

Set-TextCell $ws.Range("D48") "113.17"
Set-TextCell $ws.Range("E48") "  -3.90%  "

Set-TextCell $ws.Range("E49") "  -1.78%  "

Set-TextCell $ws.Range("D50") "5.28"
Set-TextCell $ws.Range("E50") "  -1.43%  "

Set-TextCell $ws.Range("D51") "1.602.66"
Set-TextCell $ws.Range("E51") "  +2.90%  "
